$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4").Value = "2016-08-31 18:52:21"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H4").Value = "2016-08-31 18:52:16"
$zhcn.Range("K4").Value = "2016-08-31 18:52:35"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H4").Value = "2016-08-31 18:52:21"
$dede.Range("K4").Value = "2016-08-31 18:52:42"
